# Update the "Dt. Referencia" (column G) for every data row from
# 2024-10-16 to 2024-10-17, and refresh a handful of balance values
# that were recalculated for the new reference date, then rename the
# worksheet tab to match the new export timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = Get-Date -Year 2024 -Month 10 -Day 17 -Hour 0 -Minute 0 -Second 0

$lastRow = 274
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 7).Value = $newDate.Date
}

# Row 15 - PAULO FERNANDO ULIAN: Saldo Previsto / Vl. Total updated
$ws.Cells.Item(15, 5).Value = 110286.22
$ws.Cells.Item(15, 8).Value = 110286.22

# Row 104 - ALEXANDRE FUCKNER ARTIAGA: Vl. Projetado / Saldo Previsto / Vl. Total updated
$ws.Cells.Item(104, 4).Value = 0
$ws.Cells.Item(104, 5).Value = 3769.5
$ws.Cells.Item(104, 8).Value = 3769.5

# Row 189 - FABIO HENRIQUE GOLAS: Saldo Previsto / Vl. Total updated
$ws.Cells.Item(189, 5).Value = 0
$ws.Cells.Item(189, 8).Value = 0

# Row 224 - MAURICIO ANTONIO LOPES: Saldo Previsto / Vl. Total updated
$ws.Cells.Item(224, 5).Value = 643.46
$ws.Cells.Item(224, 8).Value = 643.46

# Rename the worksheet to reflect the new export run identifier
$ws.Name = "IClientBalance-20241017-090128-"
